$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 14 de Abril de 2020 a las 11:22"

# --- España (row 5) ---
$ws.Cells.Item(5, 2).Value = 172541
$ws.Cells.Item(5, 3).Value = 2442
$ws.Cells.Item(5, 4).Value = 67504
$ws.Cells.Item(5, 5).Value = 86981
$ws.Cells.Item(5, 6).Value = 7371
$ws.Cells.Item(5, 7).Value = 300
$ws.Cells.Item(5, 8).Value = 18056

# --- Belgica (row 13) ---
$ws.Cells.Item(13, 2).Value = 31119
$ws.Cells.Item(13, 3).Value = 530
$ws.Cells.Item(13, 4).Value = 6868
$ws.Cells.Item(13, 5).Value = 20094
$ws.Cells.Item(13, 6).Value = 1223
$ws.Cells.Item(13, 7).Value = 254
$ws.Cells.Item(13, 8).Value = 4157

# --- Austria (row 20) ---
$ws.Cells.Item(20, 2).Value = 14106
$ws.Cells.Item(20, 3).Value = 65
$ws.Cells.Item(20, 5).Value = 6089
$ws.Cells.Item(20, 6).Value = 243

# --- Rows 39-41: Malasia overtakes Arabia Saudita in the ranking ---
# Row 39 becomes Malasia (was Arabia Saudita)
$ws.Cells.Item(39, 1).Value = "Malasia"
$ws.Cells.Item(39, 2).Value = 4987
$ws.Cells.Item(39, 3).Value = 170
$ws.Cells.Item(39, 4).Value = 2478
$ws.Cells.Item(39, 5).Value = 2427
$ws.Cells.Item(39, 6).Value = 60
$ws.Cells.Item(39, 7).Value = 5
$ws.Cells.Item(39, 8).Value = 82

# Row 40 becomes Arabia Saudita (was Malasia)
$ws.Cells.Item(40, 1).Value = "Arabia Saudita"
$ws.Cells.Item(40, 2).Value = 4934
$ws.Cells.Item(40, 3).Value = 0
$ws.Cells.Item(40, 4).Value = 805
$ws.Cells.Item(40, 5).Value = 4064
$ws.Cells.Item(40, 6).Value = 59
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 65

# Row 41 Indonesia (unchanged name, updated stats)
$ws.Cells.Item(41, 2).Value = 4839
$ws.Cells.Item(41, 3).Value = 282
$ws.Cells.Item(41, 4).Value = 426
$ws.Cells.Item(41, 5).Value = 3954
$ws.Cells.Item(41, 7).Value = 60
$ws.Cells.Item(41, 8).Value = 459

# --- Finlandia (row 50) ---
$ws.Cells.Item(50, 2).Value = 3161
$ws.Cells.Item(50, 3).Value = 97
$ws.Cells.Item(50, 5).Value = 2802

# --- Rows 76-77: Hong Kong overtakes Banglades in the ranking ---
# Row 76 becomes Hong Kong (was Banglades)
$ws.Cells.Item(76, 1).Value = "Hong Kong"
$ws.Cells.Item(76, 2).Value = 1013
$ws.Cells.Item(76, 3).Value = 3
$ws.Cells.Item(76, 4).Value = 434
$ws.Cells.Item(76, 5).Value = 575
$ws.Cells.Item(76, 6).Value = 13
$ws.Cells.Item(76, 7).Value = 0
$ws.Cells.Item(76, 8).Value = 4

# Row 77 becomes Banglades (was Hong Kong)
$ws.Cells.Item(77, 1).Value = "Banglades"
$ws.Cells.Item(77, 2).Value = 1012
$ws.Cells.Item(77, 3).Value = 209
$ws.Cells.Item(77, 4).Value = 42
$ws.Cells.Item(77, 5).Value = 924
$ws.Cells.Item(77, 6).Value = 1
$ws.Cells.Item(77, 7).Value = 7
$ws.Cells.Item(77, 8).Value = 46
